# Author's commit swapped the order of the two tokens in the
# "Recorded By" (column G) attribution string, turning
#   "System, dnasr281@gmail.com"
# into
#   "dnasr281@gmail.com, System"
# across every row of the session-analysis sheet that had the old
# wording (some rows already only show a single name/"System" and are
# left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "System, dnasr281@gmail.com"
$newText = "dnasr281@gmail.com, System"

$lastRow = $ws.UsedRange.Rows.Count
$colG = 7

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $colG)
    if ($cell.Value2 -eq $oldText) {
        $cell.Value = $newText
    }
}
